$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - row 3 (小野爷爷&娃展2.0), row 5 (良牙动漫秋季盛典), row 6 (熊喵M动漫嘉年华)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 50
$wsExhibition.Range("F5").Value = 3970
$wsExhibition.Range("F6").Value = 34

# Sheet "全部类型" (All types) - row 3, row 5, row 8 (same events as above, combined sheet)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 50
$wsAll.Range("F5").Value = 3970
$wsAll.Range("F8").Value = 34

$wb.Save()
